# Add Item Config And So on
# Adds a new "Icon" property row (row 19) to the Property sheet, matching
# the style/format of the existing rows (e.g. row 18 / "Extend").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hidden helper name left behind by the MySQL-for-Excel add-in (present in
# the saved workbook alongside the data edit).
$mysqlDateName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$mysqlDateName.Visible = $false

$row = 19

$ws.Cells.Item($row, 1).Value = "Icon"
$ws.Cells.Item($row, 2).Value = "string"
$ws.Cells.Item($row, 3).Value = $false
$ws.Cells.Item($row, 4).Value = $false
$ws.Cells.Item($row, 5).Value = $false
$ws.Cells.Item($row, 6).Value = $true
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = "Friend"
$ws.Cells.Item($row, 10).Value = "物品显示Icon"

# Match the text-formatted columns (A, B, I, J) used throughout the sheet.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 10).NumberFormat = "@"

# Reflect the selection left behind by the edit (matches the saved file).
[void]$ws.Range("J17").Select()
